# إضافة حدث جديد في Card18
# Row 30: several previously-empty cells get an explicit "nan" text value.
# Row 31 (new row): a brand new service-log entry is appended.
#
# Note: plain Range.Value assignment of a numeric-looking string (e.g. "18")
# or an empty string gets auto-coerced by the COM layer (exactly like typing
# into real Excel) into a Number cell / a cleared cell. The source workbook
# stores every one of these values as literal TEXT, so each write below is
# done with a leading apostrophe (Excel's "quote prefix", forcing text
# interpretation) and immediately followed by ClearFormats() so the cell
# doesn't keep the transient quote-prefix / text-number-format style that
# assigning the value would otherwise leave behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

function Set-TextValue {
    param($range, $value)
    $r = $ws.Range($range)
    $r.Value = "'" + $value
    $r.ClearFormats()
}

# ---- Row 30: fill the previously-blank cells with the literal text "nan" ----
Set-TextValue "B30" "nan"
Set-TextValue "C30" "nan"
Set-TextValue "D30" "nan"
Set-TextValue "E30" "nan"
Set-TextValue "F30" "nan"
Set-TextValue "G30" "nan"
Set-TextValue "H30" "nan"
Set-TextValue "I30" "nan"
Set-TextValue "J30" "nan"
Set-TextValue "K30" "nan"
Set-TextValue "N30" "nan"
Set-TextValue "Q30" "nan"

# ---- Row 31: new event entry ----
Set-TextValue "A31" "18"
Set-TextValue "B31" ""
Set-TextValue "C31" ""
Set-TextValue "D31" ""
Set-TextValue "E31" ""
Set-TextValue "F31" ""
Set-TextValue "G31" ""
Set-TextValue "H31" ""
Set-TextValue "I31" ""
Set-TextValue "J31" ""
Set-TextValue "K31" ""
Set-TextValue "L31" "31/1/2026"
Set-TextValue "M31" "انقطاع سير دوفر 1200"
Set-TextValue "N31" ""
Set-TextValue "O31" "تم تغير سير دوفر 1200"
Set-TextValue "P31" "عمر"
Set-TextValue "Q31" ""
